# "Generate Report for Handoff"
# Adds a new handed-off file (88686c6e-...) as a new data row on all three
# sheets (Overview, zh-cn, de-de), extending each sheet's table by one row.

$wb = $excel.ActiveWorkbook

$guidFile        = "88686c6e-a824-47d7-ab92-3c81fc031a63ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$guidFileDisplay = "e2e\88686c6e-a824-47d7-ab92-3c81fc031a63ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"
$guidFileUrl     = "https://github.com/OpenLocalizationTestOrg/oltest/blob/937f3fd158af4ba7aed3a7e91d9c1740a9db0b4f/e2e/88686c6e-a824-47d7-ab92-3c81fc031a63ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md"

$status        = "Ready for handoff"
$overviewDate  = "2016-08-12 22:34:08"
$zhHandoffDate = "2016-08-12 22:33:56"
$deHandoffDate = "2016-08-12 22:34:08"
$zhXlf         = "88686c6e-a824-47d7-ab92-3c81fc031a63ooooooooooooooooooooooooooooooooooooooooooo.5270ec627881f848d7587e3cce2d1b9467d53226.zh-cn.xlf"
$deXlf         = "88686c6e-a824-47d7-ab92-3c81fc031a63ooooooooooooooooooooooooooooooooooooooooooo.5270ec627881f848d7587e3cce2d1b9467d53226.de-de.xlf"
$zeroDate      = "0001-01-01 00:00:00"

$linkColor = 15570276  # BGR-encoded 0x6495ED to match the existing hyperlink font color

function Set-LinkLook($cell) {
    $cell.Font.Underline = $true
    $cell.Font.Color = $linkColor
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$tblOverview = $wsOverview.ListObjects.Item(1)
$tblOverview.ListRows.Add() | Out-Null

$wsOverview.Cells.Item(3, 1).Value = $guidFile
$wsOverview.Cells.Item(3, 2).Value = $guidFileDisplay
$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item(3, 2), $guidFileUrl, "", "", $guidFileDisplay) | Out-Null
Set-LinkLook $wsOverview.Cells.Item(3, 2)
$wsOverview.Cells.Item(3, 3).Value = ".md"
$wsOverview.Cells.Item(3, 4).Value = ""
$wsOverview.Cells.Item(3, 5).Value = $status
$wsOverview.Cells.Item(3, 6).Value = $status
$wsOverview.Cells.Item(3, 7).Value = $overviewDate
$wsOverview.Cells.Item(3, 7).NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Columns.Item(5).ColumnWidth = 16.38
$wsOverview.Columns.Item(6).ColumnWidth = 16.38

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$tblZh = $wsZh.ListObjects.Item(1)
$tblZh.ListRows.Add() | Out-Null

$wsZh.Cells.Item(3, 1).Value = $guidFile
$wsZh.Hyperlinks.Add($wsZh.Cells.Item(3, 1), $guidFileUrl, "", "", $guidFile) | Out-Null
Set-LinkLook $wsZh.Cells.Item(3, 1)
$wsZh.Cells.Item(3, 2).Value = ".md"
$wsZh.Cells.Item(3, 3).Value = $status
$wsZh.Cells.Item(3, 4).Value = "e2e"
$wsZh.Cells.Item(3, 5).Value = "ht"
$wsZh.Cells.Item(3, 6).Value = "False"
$wsZh.Cells.Item(3, 7).Value = $zhXlf
$wsZh.Cells.Item(3, 8).Value = $zhHandoffDate
$wsZh.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(3, 9).Value = ""
$wsZh.Cells.Item(3, 10).Value = ""
$wsZh.Cells.Item(3, 11).Value = $zeroDate
$wsZh.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(3, 12).Value = ""
$wsZh.Cells.Item(3, 13).Value = "True"
$wsZh.Cells.Item(3, 14).Value = ""
$wsZh.Cells.Item(3, 15).Value = "False"
$wsZh.Cells.Item(3, 16).Value = ""

$wsZh.Columns.Item(3).ColumnWidth = 16.38

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$tblDe = $wsDe.ListObjects.Item(1)
$tblDe.ListRows.Add() | Out-Null

$wsDe.Cells.Item(3, 1).Value = $guidFile
$wsDe.Hyperlinks.Add($wsDe.Cells.Item(3, 1), $guidFileUrl, "", "", $guidFile) | Out-Null
Set-LinkLook $wsDe.Cells.Item(3, 1)
$wsDe.Cells.Item(3, 2).Value = ".md"
$wsDe.Cells.Item(3, 3).Value = $status
$wsDe.Cells.Item(3, 4).Value = "e2e"
$wsDe.Cells.Item(3, 5).Value = "ht"
$wsDe.Cells.Item(3, 6).Value = "False"
$wsDe.Cells.Item(3, 7).Value = $deXlf
$wsDe.Cells.Item(3, 8).Value = $deHandoffDate
$wsDe.Cells.Item(3, 8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(3, 9).Value = ""
$wsDe.Cells.Item(3, 10).Value = ""
$wsDe.Cells.Item(3, 11).Value = $zeroDate
$wsDe.Cells.Item(3, 11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(3, 12).Value = ""
$wsDe.Cells.Item(3, 13).Value = "True"
$wsDe.Cells.Item(3, 14).Value = ""
$wsDe.Cells.Item(3, 15).Value = "False"
$wsDe.Cells.Item(3, 16).Value = ""

$wsDe.Columns.Item(3).ColumnWidth = 16.38
